$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 133996
$ws.Cells.Item(2,2).Value = 8
$ws.Cells.Item(2,3).Value = 16
$ws.Cells.Item(2,4).Value = 0.007395818429999999
$ws.Cells.Item(2,5).Value = 0.008847543613465568
$ws.Cells.Item(2,6).Value = 194.5253676836523
$ws.Cells.Item(2,7).Value = 0.1962899978150989
$ws.Cells.Item(2,8).Value = 43829
$ws.Cells.Item(2,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2,9).Value = 43885
$ws.Cells.Item(2,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2,10).Value = 56
$ws.Cells.Item(2,10).NumberFormat = "0"

$ws.Cells.Item(3,1).Value = 168326
$ws.Cells.Item(3,2).Value = 21
$ws.Cells.Item(3,3).Value = 43
$ws.Cells.Item(3,4).Value = 0.0058863805
$ws.Cells.Item(3,5).Value = 0.01129966615247247
$ws.Cells.Item(3,6).Value = 911.1967207380802
$ws.Cells.Item(3,7).Value = 0.9196289048036337
$ws.Cells.Item(3,8).Value = 43920
$ws.Cells.Item(3,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3,9).Value = 44074
$ws.Cells.Item(3,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3,10).Value = 154
$ws.Cells.Item(3,10).NumberFormat = "0"

$ws.Cells.Item(4,1).Value = 91885
$ws.Cells.Item(4,2).Value = 47
$ws.Cells.Item(4,3).Value = 75
$ws.Cells.Item(4,4).Value = 0.01078503426
$ws.Cells.Item(4,5).Value = 0.05716871257704084
$ws.Cells.Item(4,6).Value = 4261.964282161297
$ws.Cells.Item(4,7).Value = 4.300744642886358
$ws.Cells.Item(4,8).Value = 44102
$ws.Cells.Item(4,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4,9).Value = 44298
$ws.Cells.Item(4,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4,10).Value = 196
$ws.Cells.Item(4,10).NumberFormat = "0"

$ws.Cells.Item(5,1).Value = 90657
$ws.Cells.Item(5,2).Value = 46
$ws.Cells.Item(5,3).Value = 75
$ws.Cells.Item(5,4).Value = 0.01093120028
$ws.Cells.Item(5,5).Value = 0.05716871257704084
$ws.Cells.Item(5,6).Value = 4191.754152312831
$ws.Cells.Item(5,7).Value = 4.229865990255266
$ws.Cells.Item(5,8).Value = 44095
$ws.Cells.Item(5,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5,9).Value = 44298
$ws.Cells.Item(5,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5,10).Value = 203
$ws.Cells.Item(5,10).NumberFormat = "0"

$ws.Cells.Item(6,1).Value = 24846
$ws.Cells.Item(6,2).Value = 91
$ws.Cells.Item(6,3).Value = 96
$ws.Cells.Item(6,4).Value = 0.03989012027
$ws.Cells.Item(6,5).Value = 0.04546379302257043
$ws.Cells.Item(6,6).Value = 138.4834732103649
$ws.Cells.Item(6,7).Value = 0.1397256442157735
$ws.Cells.Item(6,8).Value = 44410
$ws.Cells.Item(6,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6,9).Value = 44445
$ws.Cells.Item(6,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6,10).Value = 35
$ws.Cells.Item(6,10).NumberFormat = "0"

$ws.Cells.Item(7,1).Value = 27981
$ws.Cells.Item(7,2).Value = 90
$ws.Cells.Item(7,3).Value = 96
$ws.Cells.Item(7,4).Value = 0.03541640101999999
$ws.Cells.Item(7,5).Value = 0.04546379302257043
$ws.Cells.Item(7,6).Value = 281.1360756239233
$ws.Cells.Item(7,7).Value = 0.2836931961803959
$ws.Cells.Item(7,8).Value = 44403
$ws.Cells.Item(7,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7,9).Value = 44445
$ws.Cells.Item(7,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7,10).Value = 42
$ws.Cells.Item(7,10).NumberFormat = "0"

$ws.Cells.Item(8,1).Value = 25371
$ws.Cells.Item(8,2).Value = 84
$ws.Cells.Item(8,3).Value = 96
$ws.Cells.Item(8,4).Value = 0.03905958055999999
$ws.Cells.Item(8,5).Value = 0.04546379302257043
$ws.Cells.Item(8,6).Value = 162.4812743878744
$ws.Cells.Item(8,7).Value = 0.1639600930361458
$ws.Cells.Item(8,8).Value = 44361
$ws.Cells.Item(8,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8,9).Value = 44445
$ws.Cells.Item(8,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8,10).Value = 84
$ws.Cells.Item(8,10).NumberFormat = "0"

$ws.Cells.Item(9,1).Value = 27656
$ws.Cells.Item(9,2).Value = 83
$ws.Cells.Item(9,3).Value = 96
$ws.Cells.Item(9,4).Value = 0.03583210630999999
$ws.Cells.Item(9,5).Value = 0.04546379302257043
$ws.Cells.Item(9,6).Value = 266.3739277228479
$ws.Cells.Item(9,7).Value = 0.2688004614979174
$ws.Cells.Item(9,8).Value = 44354
$ws.Cells.Item(9,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9,9).Value = 44445
$ws.Cells.Item(9,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9,10).Value = 91
$ws.Cells.Item(9,10).NumberFormat = "0"

$ws.Cells.Item(10,1).Value = 18112
$ws.Cells.Item(10,2).Value = 101
$ws.Cells.Item(10,3).Value = 106
$ws.Cells.Item(10,4).Value = 0.05471366900999999
$ws.Cells.Item(10,5).Value = 0.06169200115860862
$ws.Cells.Item(10,6).Value = 126.3915518755994
$ws.Cells.Item(10,7).Value = 0.1275427562230051
$ws.Cells.Item(10,8).Value = 44480
$ws.Cells.Item(10,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10,9).Value = 44515
$ws.Cells.Item(10,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10,10).Value = 35
$ws.Cells.Item(10,10).NumberFormat = "0"

$ws.Cells.Item(11,1).Value = 20539
$ws.Cells.Item(11,2).Value = 100
$ws.Cells.Item(11,3).Value = 106
$ws.Cells.Item(11,4).Value = 0.04824821000999999
$ws.Cells.Item(11,5).Value = 0.06169200115860862
$ws.Cells.Item(11,6).Value = 276.1220264012725
$ws.Cells.Item(11,7).Value = 0.2786381328099476
$ws.Cells.Item(11,8).Value = 44473
$ws.Cells.Item(11,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11,9).Value = 44515
$ws.Cells.Item(11,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11,10).Value = 42
$ws.Cells.Item(11,10).NumberFormat = "0"

$ws.Cells.Item(12,1).Value = 19127
$ws.Cells.Item(12,2).Value = 96
$ws.Cells.Item(12,3).Value = 106
$ws.Cells.Item(12,4).Value = 0.05180863687999999
$ws.Cells.Item(12,5).Value = 0.06169200115860862
$ws.Cells.Item(12,6).Value = 189.0391085569471
$ws.Cells.Item(12,7).Value = 0.1907667306804584
$ws.Cells.Item(12,8).Value = 44445
$ws.Cells.Item(12,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12,9).Value = 44515
$ws.Cells.Item(12,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12,10).Value = 70
$ws.Cells.Item(12,10).NumberFormat = "0"

$ws.Cells.Item(13,1).Value = 20105
$ws.Cells.Item(13,2).Value = 94
$ws.Cells.Item(13,3).Value = 106
$ws.Cells.Item(13,4).Value = 0.04928845922
$ws.Cells.Item(13,5).Value = 0.06169200115860862
$ws.Cells.Item(13,6).Value = 249.3732106757263
$ws.Cells.Item(13,7).Value = 0.2516520527299335
$ws.Cells.Item(13,8).Value = 44431
$ws.Cells.Item(13,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13,9).Value = 44515
$ws.Cells.Item(13,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13,10).Value = 84
$ws.Cells.Item(13,10).NumberFormat = "0"

$ws.Cells.Item(14,1).Value = 21075
$ws.Cells.Item(14,2).Value = 93
$ws.Cells.Item(14,3).Value = 106
$ws.Cells.Item(14,4).Value = 0.04702079382
$ws.Cells.Item(14,5).Value = 0.06169200115860862
$ws.Cells.Item(14,6).Value = 309.1956946611766
$ws.Cells.Item(14,7).Value = 0.3120153052875154
$ws.Cells.Item(14,8).Value = 44424
$ws.Cells.Item(14,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14,9).Value = 44515
$ws.Cells.Item(14,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14,10).Value = 91
$ws.Cells.Item(14,10).NumberFormat = "0"

$ws.Cells.Item(15,1).Value = 22605
$ws.Cells.Item(15,2).Value = 92
$ws.Cells.Item(15,3).Value = 106
$ws.Cells.Item(15,4).Value = 0.04383815435999999
$ws.Cells.Item(15,5).Value = 0.06169200115860862
$ws.Cells.Item(15,6).Value = 403.5862068825479
$ws.Cells.Item(15,7).Value = 0.4072673008081591
$ws.Cells.Item(15,8).Value = 44417
$ws.Cells.Item(15,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15,9).Value = 44515
$ws.Cells.Item(15,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15,10).Value = 98
$ws.Cells.Item(15,10).NumberFormat = "0"

$ws.Cells.Item(16,1).Value = 17497
$ws.Cells.Item(16,2).Value = 78
$ws.Cells.Item(16,3).Value = 106
$ws.Cells.Item(16,4).Value = 0.05663478820999999
$ws.Cells.Item(16,5).Value = 0.06169200115860862
$ws.Cells.Item(16,6).Value = 88.48605496180515
$ws.Cells.Item(16,7).Value = 0.08929516836642248
$ws.Cells.Item(16,8).Value = 44319
$ws.Cells.Item(16,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16,9).Value = 44515
$ws.Cells.Item(16,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16,10).Value = 196
$ws.Cells.Item(16,10).NumberFormat = "0"

$ws.Cells.Item(17,1).Value = 20936
$ws.Cells.Item(17,2).Value = 113
$ws.Cells.Item(17,3).Value = 115
$ws.Cells.Item(17,4).Value = 0.04733346617999999
$ws.Cells.Item(17,5).Value = 0.04307166
$ws.Cells.Item(17,6).Value = -89.22517418447973
$ws.Cells.Item(17,7).Value = -0.09003790603023165
$ws.Cells.Item(17,8).Value = 44564
$ws.Cells.Item(17,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17,9).Value = 44578
$ws.Cells.Item(17,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17,10).Value = 14
$ws.Cells.Item(17,10).NumberFormat = "0"

$ws.Cells.Item(18,1).Value = 19497
$ws.Cells.Item(18,2).Value = 112
$ws.Cells.Item(18,3).Value = 115
$ws.Cells.Item(18,4).Value = 0.05082625548
$ws.Cells.Item(18,5).Value = 0.04307166
$ws.Cells.Item(18,6).Value = -151.1913480735599
$ws.Cells.Item(18,7).Value = -0.1525706626775094
$ws.Cells.Item(18,8).Value = 44557
$ws.Cells.Item(18,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18,9).Value = 44578
$ws.Cells.Item(18,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18,10).Value = 21
$ws.Cells.Item(18,10).NumberFormat = "0"

$ws.Cells.Item(19,1).Value = 15110
$ws.Cells.Item(19,2).Value = 106
$ws.Cells.Item(19,3).Value = 115
$ws.Cells.Item(19,4).Value = 0.06558462911
$ws.Cells.Item(19,5).Value = 0.04307166
$ws.Cells.Item(19,6).Value = -340.1709632521
$ws.Cells.Item(19,7).Value = -0.3432659361729521
$ws.Cells.Item(19,8).Value = 44515
$ws.Cells.Item(19,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19,9).Value = 44578
$ws.Cells.Item(19,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19,10).Value = 63
$ws.Cells.Item(19,10).NumberFormat = "0"

$ws.Cells.Item(20,1).Value = 16268
$ws.Cells.Item(20,2).Value = 103
$ws.Cells.Item(20,3).Value = 115
$ws.Cells.Item(20,4).Value = 0.06091307221999999
$ws.Cells.Item(20,5).Value = 0.04307166
$ws.Cells.Item(20,6).Value = -290.2440939949598
$ws.Cells.Item(20,7).Value = -0.2928995627664631
$ws.Cells.Item(20,8).Value = 44494
$ws.Cells.Item(20,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20,9).Value = 44578
$ws.Cells.Item(20,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(20,10).Value = 84
$ws.Cells.Item(20,10).NumberFormat = "0"

$ws.Cells.Item(21,1).Value = 16090
$ws.Cells.Item(21,2).Value = 102
$ws.Cells.Item(21,3).Value = 115
$ws.Cells.Item(21,4).Value = 0.06158984831999999
$ws.Cells.Item(21,5).Value = 0.04307166
$ws.Cells.Item(21,6).Value = -297.9576500687998
$ws.Cells.Item(21,7).Value = -0.3006694905917896
$ws.Cells.Item(21,8).Value = 44487
$ws.Cells.Item(21,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21,9).Value = 44578
$ws.Cells.Item(21,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21,10).Value = 91
$ws.Cells.Item(21,10).NumberFormat = "0"

$ws.Cells.Item(22,1).Value = 16998
$ws.Cells.Item(22,2).Value = 79
$ws.Cells.Item(22,3).Value = 115
$ws.Cells.Item(22,4).Value = 0.05829907082999999
$ws.Cells.Item(22,5).Value = 0.04307166
$ws.Cells.Item(22,6).Value = -258.8355292883398
$ws.Cells.Item(22,7).Value = -0.2611947431272635
$ws.Cells.Item(22,8).Value = 44326
$ws.Cells.Item(22,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22,9).Value = 44578
$ws.Cells.Item(22,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22,10).Value = 252
$ws.Cells.Item(22,10).NumberFormat = "0"

$ws.Cells.Item(23,1).Value = 16499
$ws.Cells.Item(23,2).Value = 75
$ws.Cells.Item(23,3).Value = 115
$ws.Cells.Item(23,4).Value = 0.0600587988
$ws.Cells.Item(23,5).Value = 0.04307166
$ws.Cells.Item(23,6).Value = -280.2708030612
$ws.Cells.Item(23,7).Value = -0.2828418006921576
$ws.Cells.Item(23,8).Value = 44298
$ws.Cells.Item(23,8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23,9).Value = 44578
$ws.Cells.Item(23,9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23,10).Value = 280
$ws.Cells.Item(23,10).NumberFormat = "0"

Write-Host "done"